# Update the cryptos worksheet with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text would otherwise be auto-coerced into a
# number by Excel need to keep a text number format so the stored value
# stays an exact string (matching the source feed's formatting).
$textPriceCells = @("D5","D6","D11","D14","D20","D21","D30","D31","D32","D33","D34","D35","D38","D39","D41","D42","D43","D47")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.814.56"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.751.27"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.27%  "

# Row 5 - BNB
$ws.Range("D5").Value = "573.16"
$ws.Range("E5").Value = "  -1.52%  "

# Row 6 - Solana
$ws.Range("D6").Value = "157.41"
$ws.Range("E6").Value = "  -0.27%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.26%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.74%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -3.48%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  +0.23%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  -16.40%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -2.97%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.236.77"
$ws.Range("E13").Value = "  +0.01%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "26.45"
$ws.Range("E14").Value = "  -2.14%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "63.471.88"
$ws.Range("E15").Value = "  -1.06%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -2.78%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.752.65"
$ws.Range("E17").Value = "  -0.68%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +0.54%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -2.68%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "354.64"
$ws.Range("E20").Value = "  -2.40%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.73"
$ws.Range("E21").Value = "  -4.32%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +0.65%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  -2.48%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -1.16%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("E27").Value = "  -2.59%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  -1.09%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -4.27%  "

# Row 30 - Aptos
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  -3.01%  "

# Row 31 - Monero
$ws.Range("D31").Value = "168.98"
$ws.Range("E31").Value = "  -2.91%  "

# Row 32 - Fetch.AI
$ws.Range("D32").Value = "1.20"
$ws.Range("E32").Value = "  -6.82%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "20.13"
$ws.Range("E33").Value = "  -2.51%  "

# Row 34 - now NEARProtocol (was USDe)
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "4.86"
$ws.Range("E34").Value = "  -1.07%  "

# Row 35 - now USDe (was NEARProtocol)
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.19%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -1.23%  "

# Row 37 - Stacks
$ws.Range("E37").Value = "  -2.49%  "

# Row 38 - SuiNetwork
$ws.Range("D38").Value = "0.979"
$ws.Range("E38").Value = "  -4.29%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "6.16"
$ws.Range("E39").Value = "  +4.70%  "

# Row 40 - Filecoin
$ws.Range("E40").Value = "  -3.55%  "

# Row 41 - Bittensor
$ws.Range("D41").Value = "326.26"
$ws.Range("E41").Value = "  -4.96%  "

# Row 42 - OKB
$ws.Range("D42").Value = "38.86"
$ws.Range("E42").Value = "  -1.29%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "21.37"
$ws.Range("E43").Value = "  -3.37%  "

# Row 44 - Hedera
$ws.Range("E44").Value = "  -2.19%  "

# Row 45 - InjectiveProtocol
$ws.Range("E45").Value = "  -3.75%  "

# Row 46 - VeChain
$ws.Range("E46").Value = "  -1.91%  "

# Row 47 - Aave
$ws.Range("D47").Value = "135.19"
$ws.Range("E47").Value = "  -2.27%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  -4.37%  "

# Row 49 - Stellar
$ws.Range("E49").Value = "  -1.02%  "

# Row 50 - FirstDigitalUSD
$ws.Range("E50").Value = "  +0.42%  "

# Row 51 - WhiteBITCoin
$ws.Range("E51").Value = "  +0.34%  "
